$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, shifting Multi Step Index..Machine TEKS
# (old D..P) one column to the right (new E..Q).
$ws.Columns("D:D").Insert()

# New column header and width (the stored <col> width ends up 0.8333
# wider than the ColumnWidth we set here, so back off by that padding to
# land on the saved width of 35.5 used in the target file).
$ws.Range("D1").Value = "Ancillary UUID"
$ws.Columns("D:D").ColumnWidth = 34.6666666666667

# New ancillary UUID value for row 3 only.
$ws.Range("D3").Value = "e9779614-2fca-43cb-ae53-4af6d20e00ea"

# Reflect the final selected cell as recorded in the saved workbook.
$ws.Range("D3").Select()
